# Applies the cryptocurrency price/volume update described by the commit
# "Updated cryptos list on Fri Oct 27 15:28:11 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'33.936.74"
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = "'1.777.16"
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'224.72"
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'32.06"
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').Value = "'0.0702"
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').Value = "'2.033.78"
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = "'1.794.32"
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').Value = "'10.90"
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = "'33.932.12"
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = "'67.67"
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = "'242.51"
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = "'0.0₃0782"
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').Value = "'10.68"
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = "'2.06"
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('D25').Value = "'160.55"
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('D26').Value = "'16.28"
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = "'1.23"
$ws.Range('E30').Value = '  +3.52%  '
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').Value = "'3.64"
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').Value = "'1.80"
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').Value = "'1.391.72"
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  +5.44%  '
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').Value = "'0.0187"
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('E39').Value = '  +7.23%  '
$ws.Range('D40').Value = "'2.36"
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').Value = "'0.909"
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').Value = "'77.57"
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('D44').Value = "'13.27"
$ws.Range('E44').Value = '  +13.83%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = "'0.0₆0138"
$ws.Range('E45').Value = '  +14.10%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = "'1.07"
$ws.Range('E46').Value = '  +3.68%  '
$ws.Range('D47').Value = "'0.0497"
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('D48').Value = "'107.81"
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = "'1.933.48"
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('E51').Value = '  +0.50%  '
